# Flip the sign (positive -> negative) of every non-zero numeric value in the
# B2:E5 data block on every worksheet of the workbook. This corresponds to
# the commit "Changing the sign from + to - for the material recycled for
# each component".

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    for ($row = 2; $row -le 5; $row++) {
        for ($col = 2; $col -le 5; $col++) {
            $cell = $ws.Cells.Item($row, $col)
            $val = $cell.Value2
            if ($val -ne 0) {
                $cell.Value2 = -$val
            }
        }
    }
}
